# Regenerate sval data to filter save games - update B2:G5 values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    "B2" = 0.6753301551942219
    "C2" = 114.8270160096505
    "D2" = 0.8054896365839992
    "E2" = 8.660232485948974
    "G2" = 124.9680682873777

    "B3" = 0.3048080303191223
    "C3" = 0.3127903958511391
    "D3" = 0.1575252929769615
    "E3" = 8.660232485948974
    "G3" = 9.435356205096197

    "B4" = 0.127881588408715
    "C4" = 1.667794583268128
    "D4" = 9844.520545567508
    "E4" = 2367095152636972
    "G4" = 2367095152646818

    "B5" = 3.230985683306322
    "C5" = 1.667794583268128
    "D5" = 26.21740644021617
    "E5" = 0.496779210170732
    "G5" = 31.61296591696135
}

foreach ($addr in $data.Keys) {
    $ws.Range($addr).Value = $data[$addr]
}
